# "ajustando para sobrescrever a planilha"
# The sheet used to keep its data starting at row 5 (rows 2-4 were a gap)
# with three records in rows 5-7. Overwrite that: wipe everything below
# the header row and write a single new record right under the header,
# so the sheet collapses down to just header + one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old records (previously on rows 5-7) and the empty gap rows
# (2-4) above them, so nothing but the header remains below row 1.
$ws.Range("A2:D7").ClearContents()

# Write the new record directly under the header (row 2).
$ws.Range("A2").Value = "Dominic"
$ws.Range("B2").Value = "Bayer"
$ws.Range("C2").Value = "070.027.887-79"
$ws.Range("D2").Value = "Masculino"
